$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 7.1546

$ws.Range("B3").Value = 7.1224
$ws.Range("C3").Value = 7.1295

$ws.Range("C4").Value = 7.1287

$ws.Range("C5").Value = 7.1723

$ws.Range("C6").Value = 7.1732

$ws.Range("C7").Value = 7.1818

$ws.Range("C10").Value = 7.2527

$ws.Range("C11").Value = 7.2787

$ws.Range("C13").Value = 7.2818

$ws.Range("C16").Value = 7.0784

$ws.Range("C17").Value = 7.1508

$ws.Range("C18").Value = 7.2601

$ws.Range("C19").Value = 7.2559

$ws.Range("C20").Value = 7.236

$ws.Range("C21").Value = 7.2386

$ws.Range("C23").Value = 7.195

$ws.Range("C24").Value = 7.1675
